$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 17:17:57"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "87%"
$ws.Range("K2").Value = "8.8 MJ/m2"
$ws.Range("O2").Value = "0.0 °C"
$ws.Range("E3").Value = "2026-02-06 17:18:00"
$ws.Range("K3").Value = "12.5 MJ/m2"
$ws.Range("E4").Value = "2026-02-06 17:18:02"
$ws.Range("J4").Value = "997.0 hPa"
$ws.Range("O4").Value = "13.7 °C"
$ws.Range("E5").Value = "2026-02-06 17:18:05"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "67%"
$ws.Range("J5").Value = "997.3 hPa"
$ws.Range("K5").Value = "10.4 MJ/m2"
$ws.Range("O5").Value = "11.1 °C"
$ws.Range("E6").Value = "2026-02-06 17:18:07"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "47%"
$ws.Range("J6").Value = "998.5 hPa"
$ws.Range("E7").Value = "2026-02-06 17:18:10"
$ws.Range("J7").Value = "998.1 hPa"
$ws.Range("O7").Value = "11.8 °C"
$ws.Range("E8").Value = "2026-02-06 17:18:13"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "75%"
$ws.Range("K8").Value = "11.9 MJ/m2"
$ws.Range("O8").Value = "10.4 °C"
$ws.Range("E9").Value = "2026-02-06 17:18:16"
$ws.Range("O9").Value = "5.0 °C"
$ws.Range("E10").Value = "2026-02-06 17:18:18"
$ws.Range("O10").Value = "9.2 °C"
$ws.Range("E11").Value = "2026-02-06 17:18:21"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "75%"
$ws.Range("J11").Value = "998.5 hPa"
$ws.Range("E12").Value = "2026-02-06 17:18:24"
$ws.Range("K12").Value = "12.4 MJ/m2"
$ws.Range("E13").Value = "2026-02-06 17:18:26"
$ws.Range("E14").Value = "2026-02-06 17:18:29"
$ws.Range("K14").Value = "7.6 MJ/m2"
$ws.Range("O14").Value = "-4.1 °C"
$ws.Range("E15").Value = "2026-02-06 17:18:31"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "70%"
$ws.Range("J15").Value = "997.4 hPa"
$ws.Range("K15").Value = "11.7 MJ/m2"
$ws.Range("O15").Value = "10.7 °C"
$ws.Range("E16").Value = "2026-02-06 17:18:34"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "84%"
$ws.Range("O16").Value = "6.1 °C"
$ws.Range("E17").Value = "2026-02-06 17:18:36"
$ws.Range("E18").Value = "2026-02-06 17:18:39"
$ws.Range("E19").Value = "2026-02-06 17:18:42"
$ws.Range("E20").Value = "2026-02-06 17:18:44"
$ws.Range("E21").Value = "2026-02-06 17:18:47"
$ws.Range("J21").Value = "997.6 hPa"
$ws.Range("O21").Value = "8.6 °C"
$ws.Range("E22").Value = "2026-02-06 17:18:49"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "76%"
$ws.Range("K22").Value = "11.6 MJ/m2"
$ws.Range("E23").Value = "2026-02-06 17:18:51"
$ws.Range("I23").Value = "0.1 mm"
$ws.Range("J23").Value = "997.4 hPa"
$ws.Range("E24").Value = "2026-02-06 17:18:54"
$ws.Range("J24").Value = "996.8 hPa"
$ws.Range("E25").Value = "2026-02-06 17:18:57"
$ws.Range("J25").Value = "998.2 hPa"
$ws.Range("O25").Value = "4.5 °C"
$ws.Range("E26").Value = "2026-02-06 17:18:59"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "80%"
$ws.Range("E27").Value = "2026-02-06 17:19:02"
$ws.Range("J27").Value = "997.4 hPa"
$ws.Range("O27").Value = "11.0 °C"
$ws.Range("E28").Value = "2026-02-06 17:19:04"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "83%"
$ws.Range("J28").Value = "999.6 hPa"
$ws.Range("O28").Value = "4.9 °C"
$ws.Range("E29").Value = "2026-02-06 17:19:07"
$ws.Range("K29").Value = "12.2 MJ/m2"
$ws.Range("E30").Value = "2026-02-06 17:19:09"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "77%"
$ws.Range("K30").Value = "9.1 MJ/m2"
$ws.Range("E31").Value = "2026-02-06 17:19:12"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "83%"
$ws.Range("J31").Value = "998.9 hPa"
$ws.Range("O31").Value = "7.4 °C"
$ws.Range("E32").Value = "2026-02-06 17:19:14"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "50%"
$ws.Range("J32").Value = "998.8 hPa"
$ws.Range("K32").Value = "12.1 MJ/m2"
$ws.Range("E33").Value = "2026-02-06 17:19:17"
$ws.Range("O33").Value = "10.3 °C"
$ws.Range("E34").Value = "2026-02-06 17:19:19"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "74%"
$ws.Range("O34").Value = "8.8 °C"
$ws.Range("E35").Value = "2026-02-06 17:19:21"
$ws.Range("E36").Value = "2026-02-06 17:19:24"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "60%"
$ws.Range("I36").Value = "1.1 mm"
$ws.Range("J36").Value = "999.7 hPa"
$ws.Range("O36").Value = "13.3 °C"
